# NewPoll component on Fron-End
# Update the requirement bullets on the "UEX" slide (slide 2) to mark
# which user stories are implemented / pending.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$tr.Paragraphs(2,1).Runs(1,1).Text = "The user can log in and log out [x]"
$tr.Paragraphs(3,1).Runs(1,1).Text = "The user can write new polls []"
$tr.Paragraphs(4,1).Runs(1,1).Text = "Only him can change or delete his profile [x]"
$tr.Paragraphs(5,1).Runs(1,1).Text = "Only him can edit his polls content []"
$tr.Paragraphs(6,1).Runs(1,1).Text = "The user can vote any poll (this info will be stored inside the info of the poll) []"
